# working_hours.xlsx — add a new work-session row (2014-03-01, 18:15-21:30)
# before the summary block, pushing the summary rows down by one and
# extending the running totals to include it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 33; this shifts the old row 33 (blank
# separator) down to 34 and the summary rows (34-36) down to 35-37,
# and Excel automatically re-points the SUM()/ratio formulas in the
# summary rows plus the dimension to the new extent.
$ws.Rows(33).Insert()

# Populate the new data row with the same shape as the rows above it.
$ws.Cells.Item(33, 1).Value = 2014
$ws.Cells.Item(33, 2).Value = 3
$ws.Cells.Item(33, 3).Value = 1
$ws.Cells.Item(33, 4).Value = 0.76041666666666663
$ws.Cells.Item(33, 5).Value = 0.89583333333333337
$ws.Cells.Item(33, 4).NumberFormat = "hh:mm;@"
$ws.Cells.Item(33, 5).NumberFormat = "hh:mm;@"

# Extend the "time spent" / "sum in hours" formulas down into the new
# row, matching the pattern used by the block above (rows 28-32).
$ws.Cells.Item(33, 6).Formula = "=(E33-D33)*24*60"
$ws.Cells.Item(33, 7).Formula = "=F33/60"

# Match the author's final selection (one row further down, where the
# new blank separator row now sits).
$ws.Range("E34").Select()
